$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4 data: project info for "get users by project id"
$ws.Range("B4").Value = "validate admin mới có quyền mời member mới"
$ws.Range("C4").Value = "project"
$ws.Range("D4").Value = "trung bình"

# New column E header (bold, matching the other header cells' style)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "trạng thái"

# Set column D width (diff target stored width 13.28515625; closest achievable
# via the ColumnWidth->stored-width quantization of this runtime)
$ws.Columns.Item(4).ColumnWidth = 12.45

# Update selection to D3 (matches diff)
$ws.Range("D3").Select()
